# Adds two new bullet paragraphs to the "project notes" document:
#   1. "Upload file amazon s3"       - right after "More secure password
#      encryption technique" (same list, numId=1, no special run formatting).
#   2. "ability to delete account"   - right after "ability for parents to
#      use display name instead of actual name" (same list, numId=3, with
#      the paragraph-mark-only underline that the surrounding items in that
#      list already carry).
#
# Strategy: use Range.InsertParagraphAfter() to clone the anchor paragraph's
# pPr (style/numbering) cheaply and correctly, fill in the text, then replace
# that freshly-made paragraph's range with an explicit OOXML fragment via
# Range.InsertXML so the emitted run carries exactly the formatting we want
# (steers around COM auto-inheriting the insertion point's character
# formatting onto the new run).

$d = $word.ActiveDocument

function Add-BulletAfter($anchorText, $newText, $numId, $underlinePara) {
    # Locate the anchor paragraph (and its 1-based index) by its visible text.
    $anchor = $null
    $anchorIdx = 0
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text.TrimEnd("`r`a") -eq $anchorText) {
            $anchor = $p
            $anchorIdx = $i
        }
    }
    if ($anchor -eq $null) {
        throw "Anchor paragraph not found: $anchorText"
    }

    # Create a new paragraph right after it; Word clones the anchor's pPr
    # (pStyle + numPr, and any paragraph-mark rPr) onto the new, empty
    # paragraph automatically.
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($anchorIdx + 1)
    $newPara.Range.Text = $newText

    # Re-fetch it (ranges/handles can go stale after the mutation above) and
    # overwrite its OOXML outright so the run has exactly the formatting we
    # want, regardless of whatever the insertion point's "current formatting"
    # happened to be.
    $newPara2 = $d.Paragraphs.Item($anchorIdx + 1)
    $fullRange = $d.Range($newPara2.Range.Start, $newPara2.Range.End)

    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    if ($underlinePara) {
        $pPrExtra = "<w:rPr><w:u w:val=`"single`"/></w:rPr>"
    } else {
        $pPrExtra = ""
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:pPr>' +
        '<w:pStyle w:val="projectnotes"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr>' +
        $pPrExtra +
        '</w:pPr>' +
        '<w:r><w:t>' + $escaped + '</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $fullRange.InsertXML($xml)
}

Add-BulletAfter "More secure password encryption technique" "Upload file amazon s3" "1" $false
Add-BulletAfter "ability for parents to use display name instead of actual name" "ability to delete account" "3" $true

Write-Host "Final paragraph count:" $d.Paragraphs.Count
